$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SRS review")

# Set Acceptance cell E8 to "Accepted"
$ws.Range("E8").Value2 = "Accepted"

# Append additional reviewer comment to G8, preserving existing bold requirement-id runs
[string]$existing = $ws.Range("G8").Value2
$appendText = "Mina 20/02/2020: The values cannot be mentioned exactly/explicitly since they're a range. But the requirements Req_PO1_DGC_SRS_014_V01 and Req_PO1_DGC_SRS_016_V01 were updated to be more descriptive."
$newText = $existing + "`n" + $appendText
$ws.Range("G8").Value2 = $newText

# Setting Value2 resets all rich-text formatting on the cell, so we need to
# re-apply bold formatting to every occurrence of a requirement id that should be bold.
# (every Req_PO1_DGC_SRS_0xx_V01 occurrence is bold EXCEPT the "009" one referenced
# inline within Mali's quoted comment).
$cell = $ws.Range("G8")
$searchStart = 0
while ($true) {
    $idx = $newText.IndexOf("Req_PO1_DGC_SRS_", $searchStart)
    if ($idx -lt 0) { break }
    $len = 23
    $token = $newText.Substring($idx, $len)
    if ($token -ne "Req_PO1_DGC_SRS_009_V01") {
        $cell.Characters($idx + 1, $len).Font.Bold = $true
    }
    $searchStart = $idx + $len
}

# Reflect the updated selection/view state
$ws.Range("G8").Select() | Out-Null
